$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, shifting existing rows 147:170 down to 148:171.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with the new weekly data point.
$ws.Cells.Item(147, 1).Value = 4
$ws.Cells.Item(147, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(147, 3).Value = "Los Lagos"
$ws.Cells.Item(147, 4).Value = 44474
$ws.Cells.Item(147, 5).Value = 10
$ws.Cells.Item(147, 6).Value = 100112040
$ws.Cells.Item(147, 7).Value = "Cilantro"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 300
$ws.Cells.Item(147, 11).Value = 10000
$ws.Cells.Item(147, 12).Value = 10000
$ws.Cells.Item(147, 13).Value = 10000
$ws.Cells.Item(147, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(147, 15).Value = "Región Metropolitana"
$ws.Cells.Item(147, 16).Value = 278
$ws.Cells.Item(147, 17).Value = 36
$ws.Cells.Item(147, 18).Value = "Hortaliza"
